$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 431
$ws.Range("J18").Value = 499
$ws.Range("L18").Value = 499
$ws.Range("N18").Value = -1067

$ws.Range("H33").Value = 1082.5555
$ws.Range("I33").Value = 436
$ws.Range("J33").Value = 1599.8
$ws.Range("K33").Value = 436
$ws.Range("L33").Value = 1599.8
$ws.Range("M33").Value = -207
$ws.Range("N33").Value = -2057.8

$ws.Range("H40").Value = 6953
$ws.Range("J40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("N40").Value = -9350

$ws.Range("H76").Value = 7875.25
$ws.Range("I76").Value = 6333.3335
$ws.Range("J76").Value = 8389.223
$ws.Range("K76").Value = 6333.3335
$ws.Range("L76").Value = 8389.223
$ws.Range("M76").Value = -6018.3335
$ws.Range("N76").Value = -9019.223

$ws.Range("H79").Value = 7875.25
$ws.Range("I79").Value = 6333.3335
$ws.Range("J79").Value = 8389.223
$ws.Range("K79").Value = 6333.3335
$ws.Range("L79").Value = 8389.223
$ws.Range("M79").Value = -5241.3335
$ws.Range("N79").Value = -10573.223

$ws.Range("H80").Value = 2071.054
$ws.Range("I80").Value = 1395
$ws.Range("J80").Value = 2228.8
$ws.Range("K80").Value = 4185
$ws.Range("L80").Value = 6686.400000000001
$ws.Range("M80").Value = -3187
$ws.Range("N80").Value = -8682.400000000001

$ws.Range("H83").Value = 2071.054
$ws.Range("I83").Value = 1395
$ws.Range("J83").Value = 2228.8
$ws.Range("K83").Value = 12555
$ws.Range("L83").Value = 20059.2
$ws.Range("M83").Value = -7563
$ws.Range("N83").Value = -30043.2

$ws.Range("H86").Value = 2752.9473
$ws.Range("I86").Value = 2462.3635
$ws.Range("J86").Value = 3152.5
$ws.Range("K86").Value = 2462.3635
$ws.Range("L86").Value = 3152.5
$ws.Range("M86").Value = -1339.3635
$ws.Range("N86").Value = -5398.5

$ws.Range("H88").Value = 5088.25
$ws.Range("J88").Value = 5284.3335
$ws.Range("L88").Value = 5284.3335
$ws.Range("N88").Value = -6096.3335

$ws.Range("H89").Value = 2752.9473
$ws.Range("I89").Value = 2462.3635
$ws.Range("J89").Value = 3152.5
$ws.Range("K89").Value = 12311.8175
$ws.Range("L89").Value = 15762.5
$ws.Range("M89").Value = -6695.817499999999
$ws.Range("N89").Value = -26994.5

$ws.Range("H91").Value = 5088.25
$ws.Range("J91").Value = 5284.3335
$ws.Range("L91").Value = 5284.3335
$ws.Range("N91").Value = -8092.3335

$ws.Range("H100").Value = 4096.037
$ws.Range("J100").Value = 5330.0713
$ws.Range("L100").Value = 5330.0713
$ws.Range("N100").Value = -6412.0713

$ws.Range("H103").Value = 2250
$ws.Range("J103").Value = 2250
$ws.Range("L103").Value = 6750
$ws.Range("N103").Value = -7922

$ws.Range("H138").Value = 3479.8948
$ws.Range("I138").Value = 1372.4117
$ws.Range("J138").Value = 5185.952
$ws.Range("K138").Value = 4117.2351
$ws.Range("L138").Value = 15557.856
$ws.Range("M138").Value = 1022.7649
$ws.Range("N138").Value = -25837.856

$ws.Range("H141").Value = 4714.4375
$ws.Range("I141").Value = 2296.6155
$ws.Range("K141").Value = 6889.8465
$ws.Range("M141").Value = -1709.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 827.1739
$ws.Range("I97").Value = 800.2353000000001
$ws.Range("K97").Value = 800.2353000000001
$ws.Range("M97").Value = -304.2353000000001

$ws.Range("H102").Value = 1623.25
$ws.Range("I102").Value = 1415.3334
$ws.Range("J102").Value = 2247
$ws.Range("K102").Value = 1415.3334
$ws.Range("L102").Value = 2247
$ws.Range("M102").Value = 206.6666
$ws.Range("N102").Value = -5491

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3148.5
$ws.Range("I20").Value = 1151.3334
$ws.Range("J20").Value = 5145.6665
$ws.Range("K20").Value = 1151.3334
$ws.Range("L20").Value = 5145.6665
$ws.Range("M20").Value = -904.3334
$ws.Range("N20").Value = -5639.6665

$ws.Range("H94").Value = 2195.7273
$ws.Range("I94").Value = 2776.5386
$ws.Range("J94").Value = 1356.7778
$ws.Range("K94").Value = 2776.5386
$ws.Range("L94").Value = 1356.7778
$ws.Range("M94").Value = -2325.5386
$ws.Range("N94").Value = -2258.7778

$ws.Range("H134").Value = 1794.9318
$ws.Range("J134").Value = 6999
$ws.Range("L134").Value = 20997
$ws.Range("N134").Value = -26067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 21220.371
$ws.Range("J86").Value = 11803.315
$ws.Range("L86").Value = 11803.315
$ws.Range("N86").Value = -14049.315

$ws.Range("H89").Value = 21220.371
$ws.Range("J89").Value = 11803.315
$ws.Range("L89").Value = 59016.575
$ws.Range("N89").Value = -70248.57500000001

$ws.Range("H92").Value = 27120
$ws.Range("J92").Value = 27120
$ws.Range("L92").Value = 27120
$ws.Range("N92").Value = -32112

$ws.Range("H105").Value = 2753.611
$ws.Range("I105").Value = 2969.5
$ws.Range("J105").Value = 1998
$ws.Range("K105").Value = 2969.5
$ws.Range("L105").Value = 1998
$ws.Range("M105").Value = -1222.5
$ws.Range("N105").Value = -5492

$ws.Range("H132").Value = 3735.5
$ws.Range("J132").Value = 4124.25
$ws.Range("L132").Value = 12372.75
$ws.Range("N132").Value = -17432.75

$ws.Range("H133").Value = 79142.86
$ws.Range("J133").Value = 79142.86
$ws.Range("L133").Value = 79142.86
$ws.Range("N133").Value = -84202.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3303.3333
$ws.Range("I131").Value = 2053.5
$ws.Range("K131").Value = 6160.5
$ws.Range("M131").Value = -1120.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4386.5835
$ws.Range("I102").Value = 2702.5
$ws.Range("J102").Value = 6744.3
$ws.Range("K102").Value = 2702.5
$ws.Range("L102").Value = 6744.3
$ws.Range("M102").Value = -1080.5
$ws.Range("N102").Value = -9988.299999999999

$ws.Range("H113").Value = 3208.5789
$ws.Range("J113").Value = 3783.9
$ws.Range("L113").Value = 3783.9
$ws.Range("N113").Value = -8123.9

$ws.Range("H132").Value = 1399.5
$ws.Range("I132").Value = 741.3333
$ws.Range("K132").Value = 2223.9999
$ws.Range("M132").Value = 306.0001000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5238.4814
$ws.Range("I7").Value = 4765.778
$ws.Range("K7").Value = 4765.778
$ws.Range("M7").Value = -4653.778

$ws.Range("H40").Value = 9471
$ws.Range("I40").Value = 11283.818
$ws.Range("K40").Value = 11283.818
$ws.Range("M40").Value = -11147.818

$ws.Range("H61").Value = 2112.7273
$ws.Range("I61").Value = 1826.2
$ws.Range("K61").Value = 1826.2
$ws.Range("M61").Value = -1624.2

$ws.Range("H68").Value = 2959.9333
$ws.Range("I68").Value = 2742.5
$ws.Range("K68").Value = 2742.5
$ws.Range("M68").Value = -1993.5

$ws.Range("H71").Value = 2959.9333
$ws.Range("I71").Value = 2742.5
$ws.Range("K71").Value = 13712.5
$ws.Range("M71").Value = -9968.5

$ws.Range("H93").Value = 9075.933999999999
$ws.Range("I93").Value = 8709.444
$ws.Range("K93").Value = 8709.444
$ws.Range("M93").Value = -7461.444

$ws.Range("H113").Value = 2112.7273
$ws.Range("I113").Value = 1826.2
$ws.Range("K113").Value = 1826.2
$ws.Range("M113").Value = 343.8

$ws.Range("H126").Value = 5238.4814
$ws.Range("I126").Value = 4765.778
$ws.Range("K126").Value = 14297.334
$ws.Range("M126").Value = -11827.334

$ws.Range("H136").Value = 3870.8096
$ws.Range("I136").Value = 3694.3684
$ws.Range("J136").Value = 5547
$ws.Range("K136").Value = 11083.1052
$ws.Range("L136").Value = 16641
$ws.Range("M136").Value = -8533.1052
$ws.Range("N136").Value = -21741

$ws.Range("H139").Value = 875313.0600000001
$ws.Range("J139").Value = 994401.25
$ws.Range("L139").Value = 994401.25
$ws.Range("N139").Value = -1004681.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1381.6333
$ws.Range("I136").Value = 825
$ws.Range("K136").Value = 2475
$ws.Range("M136").Value = 75
